# Natmi following Dr Hou advice
# Rewrites the LR-pair result table: target clusters now include both
# M1 and M2 (previously only one target per sending cluster), and all
# dependent statistics are recomputed accordingly.

$data = @(
    @("ECs", "Ccl11", "Ccr3", "M1", 3, 1, 1.086278666666667, 3.258836, 0.01142600940532366, 0.01200515265421091, 2, 0.6666666666666666, 0.1323866666666667, 0.39716, 0.2977240407890032, 0.2977240407890032, 0.1438088117511111, 1.29427930576, 0.003401797690246116, 0.003574222558500497),
    @("ECs", "Ccl11", "Ccr3", "M2", 3, 1, 1.086278666666667, 3.258836, 0.01142600940532366, 0.01200515265421091, 3, 1, 0.3122756666666667, 0.936827, 0.7022759592109968, 0.7022759592109967, 0.3392183948191111, 3.052965553372, 0.008024211715077547, 0.008430930095710406),
    @("FAPs", "Ccl11", "Ccr3", "M1", 3, 1, 78.34131766666667, 235.023953, 0.8240322303590443, 0.8657994551309391, 2, 0.6666666666666666, 0.1323866666666667, 0.39716, 0.2977240407890032, 0.2977240407890032, 10.37134590816444, 93.34211317347999, 0.2453342053628693, 0.2577693122945005),
    @("FAPs", "Ccl11", "Ccr3", "M2", 3, 1, 78.34131766666667, 235.023953, 0.8240322303590443, 0.8657994551309391, 3, 1, 0.3122756666666667, 0.936827, 0.7022759592109968, 0.7022759592109967, 24.46408720190345, 220.176784817131, 0.578698024996175, 0.6080301428364386),
    @("M1", "Ccl11", "Ccr3", "M1", 3, 1, 0.971679, 2.915037, 0.01022059415658427, 0.01073863924961949, 2, 0.6666666666666666, 0.1323866666666667, 0.39716, 0.2977240407890032, 0.2977240407890032, 0.12863734388, 1.15773609492, 0.003042916591562744, 0.003197151069972105),
    @("M1", "Ccl11", "Ccr3", "M2", 3, 1, 0.971679, 2.915037, 0.01022059415658427, 0.01073863924961949, 3, 1, 0.3122756666666667, 0.936827, 0.7022759592109968, 0.7022759592109967, 0.303431707511, 2.730885367599, 0.007177677565021531, 0.007541488179647389),
    @("M2", "Ccl11", "Ccr3", "M1", 3, 1, 0.9124406666666668, 2.737322, 0.009597496442717394, 0.0100839589576554, 2, 0.6666666666666666, 0.1323866666666667, 0.39716, 0.2977240407890032, 0.2977240407890032, 0.1207949783911111, 1.08715480552, 0.002857405422383906, 0.003002237008023631),
    @("M2", "Ccl11", "Ccr3", "M2", 3, 1, 0.9124406666666668, 2.737322, 0.009597496442717394, 0.0100839589576554, 3, 1, 0.3122756666666667, 0.936827, 0.7022759592109968, 0.7022759592109967, 0.2849330174771111, 2.564397157294, 0.006740091020333489, 0.007081721949631771),
    @("sCs", "Ccl11", "Ccr3", "M1", 2, 1, 13.75898, 27.51796, 0.1447236696363304, 0.101372794007575, 2, 0.6666666666666666, 0.1323866666666667, 0.39716, 0.2977240407890032, 0.2977240407890032, 1.821505498933333, 10.9290329936, 0.04308771572194106, 0.03018111785800646),
    @("sCs", "Ccl11", "Ccr3", "M2", 2, 1, 13.75898, 27.51796, 0.1447236696363304, 0.101372794007575, 3, 1, 0.3122756666666667, 0.936827, 0.7022759592109968, 0.7022759592109967, 4.296594652153334, 25.77956791292, 0.1016359539143894, 0.07119167614956848)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $excelCol = $j + 1
        $ws.Cells.Item($excelRow, $excelCol).Value = $row[$j]
    }
}
